$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Refresh the panel_query_time values on the "data" sheet (re-fetch at
# 2021-10-05 14:33:19 instead of the original 10:50:18 run).
$data.Range("F2").Value = "2021-10-05 14:33:19.260006"
$data.Range("F3").Value = "2021-10-05 14:33:19.260014"
$data.Range("F4").Value = "2021-10-05 14:33:19.260017"

# Add a new "metadata" worksheet right after "data" holding panel-level
# metadata for this query.
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Branchio-oto-renal Syndrome"
$meta.Range("C2").Value = 57
$meta.Range("E2").Value = "2020-10-07T07:06:18.599697Z"
$meta.Range("F2").Value = "2021-10-05 14:33:19.256668"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/57/?format=json"

# data_version ("1.0") must land as literal text, not get coerced to the
# number 1 — stage it on a scratch cell (far outside the used range)
# formatted as Text, then paste only the value across so the destination
# keeps the sheet's default (General) number format. The helper column is
# then physically removed so it leaves no trace in the saved sheet.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "1.0"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Columns.Item(26).Delete()
$excel.CutCopyMode = $false

# Match the header/index styling used on the "data" sheet (bold, centered,
# bordered for the header row; bordered/centered for the leading index
# column) by copying formats across rather than re-deriving them.
$data.Range("B1:E1").Copy()
$meta.Range("B1:E1").PasteSpecial(-4122)
$data.Range("B1").Copy()
$meta.Range("F1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
